# Helper: write a value into a cell while forcing text storage (keeps leading
# zeros / decimal-looking strings like "24.42" from being silently converted
# to numbers). Any number-format residue this leaves behind on the cell's
# style gets cleaned up afterwards by re-applying the correct style (see the
# PasteSpecial formatting passes below), so we deliberately do NOT touch
# .Style here -- doing so would blow away styles applied earlier (e.g. the
# bold/bordered header style copied in from a template sheet).
function Set-TextCell($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (totals) sheet: insert a new top data row for 2022-Q4 and push
#    the existing quarters down by one.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Grow the table by one row; copy the format of the last existing data row
# (row 8) down into the brand-new row 9 so column A keeps its bold/boxed
# style.
$wsTotal.Range("A8").Copy()
$wsTotal.Range("A9").PasteSpecial(-4122)

# Push rows 8->9, 7->8, ..., 2->3 (write from the bottom up is not required
# since every value is a literal, not a relative formula).
$wsTotal.Range("A9").Value = 7
$wsTotal.Range("B9").Value = "2021-Q1"
$wsTotal.Range("C9").Value = 4
$wsTotal.Range("D9").Value = 0.04

$wsTotal.Range("A8").Value = 6
$wsTotal.Range("B8").Value = "2021-Q2"
$wsTotal.Range("C8").Value = 2
$wsTotal.Range("D8").Value = 0.09

$wsTotal.Range("A7").Value = 5
$wsTotal.Range("B7").Value = "2021-Q3"
$wsTotal.Range("C7").Value = 2
$wsTotal.Range("D7").Value = 0.08

$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q4"
$wsTotal.Range("C6").Value = 3
$wsTotal.Range("D6").Value = 0.08

$wsTotal.Range("A5").Value = 3
$wsTotal.Range("B5").Value = "2022-Q1"
$wsTotal.Range("C5").Value = 6
$wsTotal.Range("D5").Value = 0.29

$wsTotal.Range("A4").Value = 2
$wsTotal.Range("B4").Value = "2022-Q2"
$wsTotal.Range("C4").Value = 4
$wsTotal.Range("D4").Value = 0.35

$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B3").Value = "2022-Q3"
$wsTotal.Range("C3").Value = 5
$wsTotal.Range("D3").Value = 0.12

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 8
$wsTotal.Range("D2").Value = 0.68

# ---------------------------------------------------------------------------
# 2. Insert the brand-new "2022-Q4" sheet right after "总计" (i.e. before the
#    sheet that is currently "2022-Q3"), matching the formatting used by the
#    other quarterly detail sheets.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $wsTotal)
$newSheet.Name = "2022-Q4"

# Header row
Set-TextCell $newSheet "B1" "基金代码"
Set-TextCell $newSheet "C1" "基金名称"
Set-TextCell $newSheet "D1" "基金规模"
Set-TextCell $newSheet "E1" "股票总仓位"
Set-TextCell $newSheet "F1" "仓位占比"
Set-TextCell $newSheet "G1" "持有市值(亿元)"
Set-TextCell $newSheet "H1" "仓位排名"

# Data rows
$newSheet.Range("A2").Value = 0
Set-TextCell $newSheet "B2" "003986"
Set-TextCell $newSheet "C2" "申万菱信中证500指数优选增强A"
Set-TextCell $newSheet "D2" "24.42"
Set-TextCell $newSheet "E2" "93.74"
Set-TextCell $newSheet "F2" "1.88"
Set-TextCell $newSheet "G2" "0.4591"
$newSheet.Range("H2").Value = 3

$newSheet.Range("A3").Value = 1
Set-TextCell $newSheet "B3" "159851"
Set-TextCell $newSheet "C3" "华宝中证金融科技主题ETF"
Set-TextCell $newSheet "D3" "1.98"
Set-TextCell $newSheet "E3" "98.27"
Set-TextCell $newSheet "F3" "3.75"
Set-TextCell $newSheet "G3" "0.0742"
$newSheet.Range("H3").Value = 5

$newSheet.Range("A4").Value = 2
Set-TextCell $newSheet "B4" "007794"
Set-TextCell $newSheet "C4" "申万菱信中证500指数优选增强C"
Set-TextCell $newSheet "D4" "2.90"
Set-TextCell $newSheet "E4" "93.74"
Set-TextCell $newSheet "F4" "1.88"
Set-TextCell $newSheet "G4" "0.0545"
$newSheet.Range("H4").Value = 3

$newSheet.Range("A5").Value = 3
Set-TextCell $newSheet "B5" "005075"
Set-TextCell $newSheet "C5" "富国研究量化精选混合"
Set-TextCell $newSheet "D5" "2.48"
Set-TextCell $newSheet "E5" "90.71"
Set-TextCell $newSheet "F5" "1.47"
Set-TextCell $newSheet "G5" "0.0365"
$newSheet.Range("H5").Value = 8

$newSheet.Range("A6").Value = 4
Set-TextCell $newSheet "B6" "516100"
Set-TextCell $newSheet "C6" "华夏中证金融科技主题ETF"
Set-TextCell $newSheet "D6" "0.60"
Set-TextCell $newSheet "E6" "97.54"
Set-TextCell $newSheet "F6" "3.72"
Set-TextCell $newSheet "G6" "0.0223"
$newSheet.Range("H6").Value = 5

$newSheet.Range("A7").Value = 5
Set-TextCell $newSheet "B7" "002872"
Set-TextCell $newSheet "C7" "华夏智胜价值成长股票C"
Set-TextCell $newSheet "D7" "2.68"
Set-TextCell $newSheet "E7" "92.97"
Set-TextCell $newSheet "F7" "0.73"
Set-TextCell $newSheet "G7" "0.0196"
$newSheet.Range("H7").Value = 8

$newSheet.Range("A8").Value = 6
Set-TextCell $newSheet "B8" "516860"
Set-TextCell $newSheet "C8" "博时中证金融科技主题ETF"
Set-TextCell $newSheet "D8" "0.27"
Set-TextCell $newSheet "E8" "98.47"
Set-TextCell $newSheet "F8" "3.76"
Set-TextCell $newSheet "G8" "0.0102"
$newSheet.Range("H8").Value = 5

$newSheet.Range("A9").Value = 7
Set-TextCell $newSheet "B9" "002871"
Set-TextCell $newSheet "C9" "华夏智胜价值成长股票A"
Set-TextCell $newSheet "D9" "0.92"
Set-TextCell $newSheet "E9" "92.97"
Set-TextCell $newSheet "F9" "0.73"
Set-TextCell $newSheet "G9" "0.0067"
$newSheet.Range("H9").Value = 8

# ---------------------------------------------------------------------------
# 3. Apply the same look-and-feel as the other quarterly sheets: bold/boxed
#    header row (B1:H1) and bold/boxed index column (A2:A9). Do this as the
#    LAST step, after every value is already in place, so the format-only
#    paste simply overlays the style without touching any of the values/
#    types set above (in particular it overwrites the transient text number
#    format used by Set-TextCell).
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(3)   # "2022-Q3" sheet -- untouched original formatting
$templateSheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2").Copy()
$newSheet.Range("A2:A9").PasteSpecial(-4122)
